$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.992.89"
$ws.Range("D3").Value = "2.049.67"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'248.45"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'57.07"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'0.0782"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "'15.82"
$ws.Range("E12").Value = "  +4.93%  "
$ws.Range("D13").Value = "2.345.97"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +7.05%  "
$ws.Range("D15").Value = "'0.794"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").Value = "2.044.59"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "36.993.28"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'16.28"
$ws.Range("E18").Value = "  +14.02%  "
$ws.Range("D19").Value = "'74.28"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "'5.34"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'236.31"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("E25").Value = "  +11.19%  "
$ws.Range("D26").Value = "'167.90"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "'9.09"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  +5.80%  "
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'0.0615"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "'0.0885"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'2.23"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "'1.74"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "'3.16"
$ws.Range("E40").Value = "  +13.35%  "
$ws.Range("D41").Value = "'4.93"
$ws.Range("E41").Value = "  +24.87%  "
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "'17.30"
$ws.Range("E43").Value = "  -5.29%  "
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "'95.74"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "'2.45"
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("D47").Value = "1.280.02"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").Value = "2.237.89"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'6.69"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'3.53"
$ws.Range("E51").Value = "  -19.86%  "
